# Rotates columns D (codeforiati:group-name), E (codeforiati:category-name) and
# F (codeforiati:category-code) one step to the left for every row in the sheet
# (including the header row):
#   new D = old F
#   new E = old D
#   new F = old E
# Column G (codeforiati:group-code) is left untouched.
#
# This reproduces the upstream codeforIATI/codelists regeneration which changed
# the column order of the SectorGroup codelist from
#   code, name, status, group-name, category-name, category-code, group-code
# to
#   code, name, status, category-code, group-name, category-name, group-code

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$colD = 4
$colE = 5
$colF = 6

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, $colD)
    $eCell = $ws.Cells.Item($r, $colE)
    $fCell = $ws.Cells.Item($r, $colF)

    $oldD = $dCell.Value2
    $oldE = $eCell.Value2
    $oldF = $fCell.Value2

    $dCell.Value2 = $oldF
    $eCell.Value2 = $oldD
    $fCell.Value2 = $oldE
}

Write-Host "Rotated columns D/E/F for $lastRow rows"
